# Insert a new daily-ranking row right before the existing "2026/12/29"
# block (old row 750), pushing rows 750:791 down to 751:792, and fill
# the newly inserted row with the 2026/01/30 entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(750).Insert()

# Column A holds dates as plain text (e.g. "2026/01/30"), not real Excel
# date serials, matching the rest of the column. Force text formatting
# before assigning so Excel doesn't auto-convert the string to a date,
# then drop back to the default "Normal" style so no stray number format
# is left attached to the cell.
$ws.Range("A750").NumberFormat = "@"
$ws.Range("A750").Value = "2026/01/30"
$ws.Range("A750").Style = "Normal"

$ws.Range("B750").Value = "金"
$ws.Range("C750").Value = 17
$ws.Range("D750").Value = 201
